{"js": "// Update the date line and the twenty-five division problems in the table\n// to the new values, per the commit's regenerated worksheet output.\nconst replacements = [\n  [\"2024-09-13 Friday\", \"2024-09-14 Saturday\"],\n  [\"533\u00f79=\", \"377\u00f79=\"],\n  [\"336\u00f79=\", \"632\u00f75=\"],\n  [\"870\u00f76=\", \"786\u00f79=\"],\n  [\"230\u00f73=\", \"206\u00f78=\"],\n  [\"248\u00f73=\", \"804\u00f73=\"],\n  [\"111\u00f77=\", \"494\u00f79=\"],\n  [\"169\u00f75=\", \"921\u00f78=\"],\n  [\"642\u00f74=\", \"135\u00f78=\"],\n  [\"275\u00f74=\", \"981\u00f79=\"],\n  [\"924\u00f74=\", \"205\u00f78=\"],\n  [\"915\u00f79=\", \"148\u00f77=\"],\n  [\"294\u00f78=\", \"655\u00f75=\"],\n  [\"408\u00f73=\", \"907\u00f77=\"],\n  [\"369\u00f78=\", \"936\u00f72=\"],\n  [\"767\u00f75=\", \"661\u00f76=\"],\n  [\"914\u00f78=\", \"931\u00f77=\"],\n  [\"667\u00f76=\", \"984\u00f75=\"],\n  [\"692\u00f75=\", \"116\u00f77=\"],\n  [\"282\u00f73=\", \"841\u00f72=\"],\n  [\"228\u00f77=\", \"609\u00f79=\"],\n  [\"856\u00f79=\", \"386\u00f74=\"],\n  [\"302\u00f77=\", \"791\u00f79=\"],\n  [\"441\u00f73=\", \"723\u00f73=\"],\n  [\"527\u00f74=\", \"654\u00f76=\"],\n  [\"300\u00f75=\", \"947\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-five division problems in the table\n# to the new values, per the commit's regenerated worksheet output.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-09-13 Friday\", \"2024-09-14 Saturday\"),\n    @(\"533\u00f79=\", \"377\u00f79=\"),\n    @(\"336\u00f79=\", \"632\u00f75=\"),\n    @(\"870\u00f76=\", \"786\u00f79=\"),\n    @(\"230\u00f73=\", \"206\u00f78=\"),\n    @(\"248\u00f73=\", \"804\u00f73=\"),\n    @(\"111\u00f77=\", \"494\u00f79=\"),\n    @(\"169\u00f75=\", \"921\u00f78=\"),\n    @(\"642\u00f74=\", \"135\u00f78=\"),\n    @(\"275\u00f74=\", \"981\u00f79=\"),\n    @(\"924\u00f74=\", \"205\u00f78=\"),\n    @(\"915\u00f79=\", \"148\u00f77=\"),\n    @(\"294\u00f78=\", \"655\u00f75=\"),\n    @(\"408\u00f73=\", \"907\u00f77=\"),\n    @(\"369\u00f78=\", \"936\u00f72=\"),\n    @(\"767\u00f75=\", \"661\u00f76=\"),\n    @(\"914\u00f78=\", \"931\u00f77=\"),\n    @(\"667\u00f76=\", \"984\u00f75=\"),\n    @(\"692\u00f75=\", \"116\u00f77=\"),\n    @(\"282\u00f73=\", \"841\u00f72=\"),\n    @(\"228\u00f77=\", \"609\u00f79=\"),\n    @(\"856\u00f79=\", \"386\u00f74=\"),\n    @(\"302\u00f77=\", \"791\u00f79=\"),\n    @(\"441\u00f73=\", \"723\u00f73=\"),\n    @(\"527\u00f74=\", \"654\u00f76=\"),\n    @(\"300\u00f75=\", \"947\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
